$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "conventies C#"
$ws.Range("B6").Value = "1 uur"

$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = 42843

$ws.Range("C7").Select()
